# Pytheas nucleotide_table.xlsx edit
# "digest map plot and massacre output support added"
#
# Adds four new reference rows (123-126) describing single-atom
# "nucleotides" (Hydrogen, Carbon, Oxygen, Nitrogen) used as digest-map /
# massacre output placeholders, and updates the active view/selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: A=Name, B=AMBER_3_letter_code, C=Pytheas_1_letter_code,
#          D=Pytheas_ID, E=Modomics_1, F=Modomics_U, G=Originating_base,
#          H=Type, I=Linkage, J=Precursors

$rows = @(
    @{ Name = "Hydrogen"; Code = "HHH" },
    @{ Name = "Carbon";   Code = "CCC" },
    @{ Name = "Oxygen";   Code = "OOO" },
    @{ Name = "Nitrogen"; Code = "NNN" }
)

$r = 123
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row.Name   # A - Name
    $ws.Cells.Item($r, 2).Value = $row.Code   # B - AMBER_3_letter_code
    $ws.Cells.Item($r, 4).Value = $row.Code   # D - Pytheas_ID
    $ws.Cells.Item($r, 8).Value = "atom"      # H - Type
    $ws.Cells.Item($r, 9).Value = "none"      # I - Linkage
    $ws.Cells.Item($r, 10).Value = "none"     # J - Precursors
    $r = $r + 1
}

# Update the active selection / scroll position to reflect where the
# editor was working when the rows were added.
$win = $excel.ActiveWindow
$win.ScrollRow = 67
$win.ScrollColumn = 1
$ws.Range("O114").Select()
